# Deploying to gh-pages: add the new "2020" data point (column N) to the
# "рус,англ" sheet, mirroring the existing 2010-2019 columns (D:M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New year header (N4), matching the formatting already used by the
# adjacent 2019 header cell (M4).
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2020

# New data value (N5), matching the formatting already used by the
# adjacent 2019 value cell (M5).
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 2.1

$excel.CutCopyMode = $false

# Leave the selection where editing finished.
$ws.Range("N9").Select()
